$d = $word.ActiveDocument

# The content to keep ends with the paragraph containing "Mostrando um
# resultado bastante considerável." followed by one blank paragraph,
# which becomes the document's new final paragraph. Every paragraph
# after that blank one (more blank paragraphs, plus the "A ideia se
# mostra...", "O estudo é de certa forma antigo..." and "Acredito que
# possa ser bastante útil..." paragraphs) is removed.

$findRange = $d.Content
$found = $findRange.Find.Execute("Mostrando um resultado bastante considerável.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor text not found"
}

$paragraphs = $d.Paragraphs
$total = $paragraphs.Count

# Work out which paragraph number (1-based) holds the anchor text by
# scanning, since Information() isn't reliable across hosts.
$anchorParaNumber = 0
for ($i = 1; $i -le $total; $i++) {
    $p = $paragraphs.Item($i)
    if ($p.Range.Text -match [regex]::Escape("Mostrando um resultado bastante considerável.")) {
        $anchorParaNumber = $i
        break
    }
}
if ($anchorParaNumber -eq 0) {
    throw "Could not locate anchor paragraph"
}

# Keep the blank paragraph immediately following the anchor paragraph;
# delete everything after it through the end of the document.
$keepParaNumber = $anchorParaNumber + 1
$keepPara = $paragraphs.Item($keepParaNumber)
$lastPara = $paragraphs.Item($total)

$deleteStart = $keepPara.Range.End
$deleteEnd = $lastPara.Range.End

if ($deleteEnd -gt $deleteStart) {
    $r = $d.Range($deleteStart, $deleteEnd)
    $r.Delete()
}
